$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (shifting old C,D to E,F)
$ws.Range("C:D").Insert()

$ws.Range("C1").Value = "Also Empty"
$ws.Range("D1").Value = "Not Empty"
$ws.Range("D2").Value = "This"
$ws.Range("D3").Value = "is not"
$ws.Range("D4").Value = "empty"

$ws.Range("J11").Select() | Out-Null
